$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H2").Value = 366.88
$ws.Range("I2").Value = 269.27777
$ws.Range("J2").Value = 617.8570999999999
$ws.Range("K2").Value = 269.27777
$ws.Range("L2").Value = 617.8570999999999
$ws.Range("M2").Value = -156.27777
$ws.Range("N2").Value = -843.8570999999999

$ws.Range("H64").Value = 6418.636
$ws.Range("I64").Value = 4681
$ws.Range("J64").Value = 7866.6665
$ws.Range("K64").Value = 4681
$ws.Range("L64").Value = 7866.6665
$ws.Range("M64").Value = -4433
$ws.Range("N64").Value = -8362.666499999999

$ws.Range("H67").Value = 6418.636
$ws.Range("I67").Value = 4681
$ws.Range("J67").Value = 7866.6665
$ws.Range("K67").Value = 4681
$ws.Range("L67").Value = 7866.6665
$ws.Range("M67").Value = -3823
$ws.Range("N67").Value = -9582.666499999999

$ws.Range("H80").Value = 2212.6
$ws.Range("I80").Value = 309.33334
$ws.Range("J80").Value = 3028.2856
$ws.Range("K80").Value = 928.0000200000001
$ws.Range("L80").Value = 9084.856800000001
$ws.Range("M80").Value = 69.99997999999994
$ws.Range("N80").Value = -11080.8568

$ws.Range("H83").Value = 2212.6
$ws.Range("I83").Value = 309.33334
$ws.Range("J83").Value = 3028.2856
$ws.Range("K83").Value = 2784.00006
$ws.Range("L83").Value = 27254.5704
$ws.Range("M83").Value = 2207.99994
$ws.Range("N83").Value = -37238.5704

$ws.Range("H120").Value = 166520.33
$ws.Range("J120").Value = 166520.33
$ws.Range("L120").Value = 166520.33
$ws.Range("N120").Value = -176196.33

$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 3031107
$ws.Range("I2").Value = 3497296.8
$ws.Range("K2").Value = 3497296.8
$ws.Range("M2").Value = -3497183.8

$ws.Range("H13").Value = 2919.4
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 2399.25
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 2399.25
$ws.Range("M13").Value = -4856
$ws.Range("N13").Value = -2687.25

$ws.Range("H63").Value = 8145.6787
$ws.Range("J63").Value = 9876.429
$ws.Range("L63").Value = 9876.429
$ws.Range("N63").Value = -11248.429

$ws.Range("H66").Value = 8145.6787
$ws.Range("J66").Value = 9876.429
$ws.Range("L66").Value = 49382.145
$ws.Range("N66").Value = -56246.145

$ws.Range("H116").Value = 3031107
$ws.Range("I116").Value = 3497296.8
$ws.Range("K116").Value = 3497296.8
$ws.Range("M116").Value = -3495002.8

$ws.Range("H117").Value = 89998.664
$ws.Range("J117").Value = 89998.664
$ws.Range("L117").Value = 89998.664
$ws.Range("N117").Value = -99176.664

$ws.Range("H138").Value = 125000
$ws.Range("J138").Value = 125000
$ws.Range("L138").Value = 125000
$ws.Range("N138").Value = -135280

$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 3031107
$ws.Range("I3").Value = 3497296.8
$ws.Range("K3").Value = 3497296.8
$ws.Range("M3").Value = -3497182.8

$ws.Range("H86").Value = 126157.06
$ws.Range("I86").Value = 1234.2
$ws.Range("J86").Value = 2000000
$ws.Range("K86").Value = 1234.2
$ws.Range("L86").Value = 2000000
$ws.Range("M86").Value = -111.2
$ws.Range("N86").Value = -2002246

$ws.Range("H89").Value = 126157.06
$ws.Range("I89").Value = 1234.2
$ws.Range("J89").Value = 2000000
$ws.Range("K89").Value = 6171
$ws.Range("L89").Value = 10000000
$ws.Range("M89").Value = -555
$ws.Range("N89").Value = -10011232

$ws.Range("H97").Value = 21459.375
$ws.Range("I97").Value = 5362.5
$ws.Range("J97").Value = 69750
$ws.Range("K97").Value = 5362.5
$ws.Range("L97").Value = 69750
$ws.Range("M97").Value = -4371.5
$ws.Range("N97").Value = -71732

$ws.Range("H106").Value = 53866.668
$ws.Range("J106").Value = 53866.668
$ws.Range("L106").Value = 53866.668
$ws.Range("N106").Value = -56390.668

$ws.Range("H111").Value = 92698.5
$ws.Range("J111").Value = 92698.5
$ws.Range("L111").Value = 92698.5
$ws.Range("N111").Value = -100878.5

$ws.Range("H117").Value = 94294.664
$ws.Range("J117").Value = 94294.664
$ws.Range("L117").Value = 94294.664
$ws.Range("N117").Value = -103472.664

$ws.Range("H120").Value = 76265
$ws.Range("J120").Value = 76265
$ws.Range("L120").Value = 76265
$ws.Range("N120").Value = -85941

$ws.Range("H125").Value = 45390
$ws.Range("J125").Value = 45390
$ws.Range("L125").Value = 45390
$ws.Range("N125").Value = -55230

$ws.Range("H134").Value = 10260.321
$ws.Range("I134").Value = 9771.559999999999
$ws.Range("K134").Value = 29314.68
$ws.Range("M134").Value = -26779.68

$ws = $wb.Worksheets.Item(4)
$ws.Range("H74").Value = 60077
$ws.Range("J74").Value = 60077
$ws.Range("L74").Value = 60077
$ws.Range("N74").Value = -61825

$ws.Range("H77").Value = 60077
$ws.Range("J77").Value = 60077
$ws.Range("L77").Value = 180231
$ws.Range("N77").Value = -188967

$ws = $wb.Worksheets.Item(5)
$ws.Range("H23").Value = 386.85715
$ws.Range("J23").Value = 781.8333
$ws.Range("L23").Value = 2345.4999
$ws.Range("N23").Value = -2815.4999

$ws.Range("H55").Value = 1939.25
$ws.Range("I55").Value = 378
$ws.Range("K55").Value = 1134
$ws.Range("M55").Value = -957

$ws.Range("H68").Value = 291571.56
$ws.Range("J68").Value = 9575
$ws.Range("L68").Value = 28725
$ws.Range("N68").Value = -30347

$ws.Range("H71").Value = 291571.56
$ws.Range("J71").Value = 9575
$ws.Range("L71").Value = 86175
$ws.Range("N71").Value = -94287

$ws = $wb.Worksheets.Item(6)
$ws.Range("H3").Value = 2500808.2
$ws.Range("I3").Value = 249
$ws.Range("J3").Value = 3000920
$ws.Range("K3").Value = 249
$ws.Range("L3").Value = 3000920
$ws.Range("M3").Value = -133
$ws.Range("N3").Value = -3001152

$ws.Range("H9").Value = 6499.5
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H10").Value = 8000
$ws.Range("J10").Value = 8000
$ws.Range("L10").Value = 8000
$ws.Range("N10").Value = -8338

$ws.Range("H11").Value = 9002923
$ws.Range("I11").Value = 669513.3
$ws.Range("J11").Value = 17336334
$ws.Range("K11").Value = 669513.3
$ws.Range("L11").Value = 17336334
$ws.Range("M11").Value = -669374.3
$ws.Range("N11").Value = -17336612

$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H13").Value = 445.66666
$ws.Range("I13").Value = 583.3333
$ws.Range("J13").Value = 376.83334
$ws.Range("K13").Value = 583.3333
$ws.Range("L13").Value = 376.83334
$ws.Range("M13").Value = -444.3333
$ws.Range("N13").Value = -654.83334

$ws.Range("H14").Value = 513405
$ws.Range("J14").Value = 17857.143
$ws.Range("L14").Value = 17857.143
$ws.Range("N14").Value = -18193.143

$ws.Range("H19").Value = 18792.8
$ws.Range("J19").Value = 18792.8
$ws.Range("L19").Value = 18792.8
$ws.Range("N19").Value = -19368.8

$ws.Range("H22").Value = 2600
$ws.Range("J22").Value = 6000.5
$ws.Range("L22").Value = 6000.5
$ws.Range("N22").Value = -7058.5

$ws.Range("H25").Value = 6500
$ws.Range("J25").Value = 6500
$ws.Range("L25").Value = 6500
$ws.Range("N25").Value = -7558

$ws.Range("H99").Value = 10745.1
$ws.Range("I99").Value = 10745.1
$ws.Range("K99").Value = 10745.1
$ws.Range("M99").Value = -8499.1

$ws.Range("H132").Value = 5341.6772
$ws.Range("I132").Value = 3195.0833
$ws.Range("K132").Value = 9585.249899999999
$ws.Range("M132").Value = -7055.249899999999

$ws.Range("H141").Value = 80000.336
$ws.Range("J141").Value = 80000.336
$ws.Range("L141").Value = 80000.336
$ws.Range("N141").Value = -90360.336

$ws = $wb.Worksheets.Item(7)
$ws.Range("H12").Value = 263.33334
$ws.Range("I12").Value = 205
$ws.Range("K12").Value = 205
$ws.Range("M12").Value = -35

$ws.Range("H124").Value = 80971
$ws.Range("J124").Value = 80971
$ws.Range("L124").Value = 80971
$ws.Range("N124").Value = -90791

$ws.Range("H125").Value = 171054.89
$ws.Range("J125").Value = 171054.89
$ws.Range("L125").Value = 171054.89
$ws.Range("N125").Value = -180894.89

$ws = $wb.Worksheets.Item(8)
$ws.Range("H19").Value = 15002.5
$ws.Range("I19").Value = 15002.5
$ws.Range("K19").Value = 15002.5
$ws.Range("M19").Value = -14828.5

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H107").Value = 907.41174
$ws.Range("I107").Value = 1143.1
$ws.Range("J107").Value = 570.7143
$ws.Range("K107").Value = 3429.3
$ws.Range("L107").Value = 1712.1429
$ws.Range("M107").Value = -1509.3
$ws.Range("N107").Value = -5552.1429

$ws.Range("H116").Value = 101489
$ws.Range("J116").Value = 101489
$ws.Range("L116").Value = 101489
$ws.Range("N116").Value = -110667

$ws.Range("H120").Value = 71419.5
$ws.Range("J120").Value = 71419.5
$ws.Range("L120").Value = 71419.5
$ws.Range("N120").Value = -81095.5

$ws.Range("H140").Value = 119000
$ws.Range("J140").Value = 119000
$ws.Range("L140").Value = 119000

$ws.Range("H141").Value = 104333.164
$ws.Range("J141").Value = 104333.164
$ws.Range("L141").Value = 104333.164
$ws.Range("N141").Value = -114693.164
